# Re-order the comma separated "Recorded By" values in column G so that
# the recognised account names sort into a fixed priority order instead
# of whatever order they happened to be recorded in:
#   0: admin@admin.com, backup@backdoor.com
#   1: dnasr281@gmail.com, system (lower-case)
#   2: System (capitalised)
# Anything not recognised keeps its relative position (treated as tier 1)
# via a stable sort. Note: names differ only by case ("system" vs
# "System"), so comparisons must be case-sensitive - hashtable lookups
# and -eq are case-insensitive in this engine, so plain .Equals() is
# used instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-NameRank($name) {
    if ($name.Equals("admin@admin.com")) { return 0 }
    if ($name.Equals("backup@backdoor.com")) { return 0 }
    if ($name.Equals("dnasr281@gmail.com")) { return 1 }
    if ($name.Equals("system")) { return 1 }
    if ($name.Equals("System")) { return 2 }
    return 1
}

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G - "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val.Split(",")
    $names = @()
    foreach ($p in $parts) { $names += $p.Trim() }

    if ($names.Count -le 1) { continue }

    $keyed = @()
    for ($i = 0; $i -lt $names.Count; $i++) {
        $name = $names[$i]
        $item = @{}
        $item["Name"] = $name
        $item["Key"] = (Get-NameRank $name) * 1000 + $i
        $keyed += $item
    }

    $sorted = $keyed | Sort-Object -Property {$_["Key"]}
    $newNames = @()
    foreach ($k in $sorted) { $newNames += $k["Name"] }
    $newVal = $newNames -join ", "

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
